$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Vector2D" class header in column Q, row 1
$ws.Range("Q1").Value = "Vector2D"

# Update selection to the newly added cell
$ws.Range("Q1").Select()
